# Commit: Sync attendance_reports, modules_schedules, and assets from main repo
#
# In the "Session Analysis Results" sheet, column G ("Recorded By") lists who
# recorded a session. Rows that were auto-recorded by the system show both
# "System" and the instructor's email, joined with ", ". This edit swaps the
# order of those two names so the instructor's email comes first:
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
# Rows recorded only by the instructor (no "System" prefix) are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$rows = 8,9,10,12,14,15,17,18,34,35,36,38,40,41,43,44,60,61,62,64,66,67,69,70,`
        86,87,88,90,92,93,95,96,112,113,114,116,118,119,121,122,138,139,140,142,`
        144,145,147,148,164,167,170,174,191,194,197,201,218,221,224,228,245,248,`
        251,255,272,275,278,282,299,302,305,309

$changed = 0
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
        $changed++
    }
}

Write-Output "Updated $changed cell(s) in column G."
